# "Generate Report for Handback"
#
# The localization report for the 63191932-... file has been handed back
# (it is now in sync with en-US), so every place that showed its status as
# "Ready for handoff" needs to flip to "Handed back: in sync with en-US",
# and the per-locale "Latest Handback DateTime" timestamps need to be
# refreshed to reflect the new handback.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 63191932-...md file; its zh-cn / de-de status
# columns move from "Ready for handoff" to "Handed back: in sync with en-US".
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn detail sheet: same status flip for the 63191932-...md row, plus the
# Latest Handback DateTime for both rows refreshes to the handback run time.
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("G2").Value = "2016-02-25 06:24:59"
$zhcn.Range("G3").Value = "2016-02-25 06:24:59"

# de-de detail sheet: same treatment.
$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("G2").Value = "2016-02-25 06:25:22"
$dede.Range("G3").Value = "2016-02-25 06:25:22"
